# "plano de ação atualizado" — fill in the SPRINT 3D rows (40-42) of the
# "Plano de Ação" sheet with the newly-defined tasks (dashboard site,
# solution diagram, presentation slides), matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 40: "Site estátido dashboard" -------------------------------------
$ws.Range("B40").Value = "Site estátido dashboard"
$ws.Range("B41").Value = "Diagrama de solução "
$ws.Range("B42").Value = "slides da apresentação"

$ws.Range("C40").Value = "Toda equipe"
$ws.Range("C41").Value = "Toda equipe"
$ws.Range("C42").Value = "Toda equipe"

$ws.Range("D40").Value = "Essencial"
$ws.Range("D41").Value = "Essencia"
$ws.Range("D42").Value = "Essencia"

$ws.Range("H40").Value = "Criação da dashboard e linkar junto ao login"
$ws.Range("H41").Value = "Terminar a criação do diagrama de solução"
$ws.Range("H42").Value = "Fazer os slides para a apresentação"

# Progress column is percentage-formatted.
$ws.Range("E40:E42").NumberFormat = "0%"
$ws.Range("E40").Value = 0.3
$ws.Range("E41").Value = 0.4
$ws.Range("E42").Value = 0

# Start / end dates for the new sprint items.
$ws.Range("F40").Value = "10/20/2023"
$ws.Range("F41").Value = "10/20/2023"
$ws.Range("F42").Value = "10/20/2023"

$ws.Range("G40").Value = "10/23/2023"
$ws.Range("G41").Value = "10/23/2023"
$ws.Range("G42").Value = "10/23/2023"

# Leave the cursor on the last-edited cell, mirroring the author's final
# selection after entering this data.
$ws.Range("H42").Select()
